$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.839.82'
$ws.Range("E2").Value = '  -1.00%  '

$ws.Range("D3").Value = '1.900.17'
$ws.Range("E3").Value = '  -0.68%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7719'
$ws.Range("E5").Value = '  +4.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.34'
$ws.Range("E6").Value = '  -1.44%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3062'
$ws.Range("E8").Value = '  -2.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.52'
$ws.Range("E9").Value = '  -5.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06859'
$ws.Range("E10").Value = '  -1.81%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07987'
$ws.Range("E11").Value = '  -0.12%  '

$ws.Range("D12").Value = '1.922.71'
$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7387'
$ws.Range("E13").Value = '  -5.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.184'
$ws.Range("E14").Value = '  -2.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.41'

$ws.Range("D16").Value = '29.858.20'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.78'
$ws.Range("E17").Value = '  -4.47%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.902'
$ws.Range("E18").Value = '  -0.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.77'
$ws.Range("E19").Value = '  +1.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007720'
$ws.Range("E20").Value = '  -1.70%  '

$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("D22").Value = '2.148.91'
$ws.Range("E22").Value = '  +0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.942'
$ws.Range("E24").Value = '  -4.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.278'
$ws.Range("E25").Value = '  -1.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.78'
$ws.Range("E26").Value = '  -0.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.75'
$ws.Range("E27").Value = '  -1.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1293'
$ws.Range("E28").Value = '  +0.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.030'
$ws.Range("E29").Value = '  -1.91%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.391'
$ws.Range("E30").Value = '  +2.72%  '

$ws.Range("E31").Value = '  -2.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.279'
$ws.Range("E32").Value = '  -1.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.067'
$ws.Range("E33").Value = '  -0.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05258'
$ws.Range("E34").Value = '  +1.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.246'
$ws.Range("E35").Value = '  -4.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7284'
$ws.Range("E36").Value = '  -3.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.726'
$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("E38").Value = '  -1.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.780'
$ws.Range("E39").Value = '  -0.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.194'
$ws.Range("E40").Value = '  -2.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4417'
$ws.Range("E41").Value = '  -2.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.12'
$ws.Range("E42").Value = '  -4.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8377'
$ws.Range("E44").Value = '  -0.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.879'
$ws.Range("E45").Value = '  -4.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.607'
$ws.Range("E46").Value = '  -3.45%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.35'
$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.764'
$ws.Range("E48").Value = '  -1.87%  '

$ws.Range("D49").Value = '2.055.35'
$ws.Range("E49").Value = '  -0.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.17'
$ws.Range("E50").Value = '  -2.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '922.10'
$ws.Range("E51").Value = '  -1.77%  '
